$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title heading + bold byline repeated later in the doc (both occurrences -> same new text)
Replace-Text "Play Marlin Catch Free Slot Game | Stakelogic" "Play Marlin Catch for Free - Exciting Underwater Adventure"

# "What we like" bullet list
Replace-Text "High volatility for a thrilling experience" "High volatility for big win potential"
Replace-Text "Sticky Wilds and multipliers help increase winnings" "Marine-themed symbols create immersive gameplay"
Replace-Text "Fun, marine-themed symbols" "Sticky Wilds and multipliers during Free Spins"
Replace-Text "Gamble and Buy Bonus options available" "Additional features like Gamble and Buy Bonus options"

# "What we don't like" bullet list
Replace-Text "Only 10 paylines may not appeal to some players" "High volatility may not be suitable for all players"
Replace-Text "Not suitable for those who prefer low volatility games" "Limited number of paylines"

# Meta description (italic run at the end)
Replace-Text "Read our review of Marlin Catch, a high volatility online slot game with features like Wilds, Scatters, and Free Spins. Play this marine-themed adventure for free!" "Read our review of Marlin Catch and play for free at top-notch ADM/AAMS casinos."

Write-Output "Done"
